$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.187.41'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '1.785.23'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '225.96'
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '31.89'
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").Value = '2.042.50'
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").Value = '11.03'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '1.786.43'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.623'
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '34.124.59'
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").Value = '68.04'
$ws.Range("E18").Value = '  +2.61%  '
$ws.Range("D19").Value = '245.95'
$ws.Range("E19").Value = '  +3.54%  '
$ws.Range("D20").Value = '0.0₃0779'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("D21").Value = '10.97'
$ws.Range("E21").Value = '  +4.12%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").Value = '161.66'
$ws.Range("E25").Value = '  +1.59%  '
$ws.Range("E26").Value = '  +2.82%  '
$ws.Range("E27").Value = '  +1.59%  '
$ws.Range("E28").Value = '  +2.01%  '
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("E31").Value = '  +1.72%  '
$ws.Range("D32").Value = '3.70'
$ws.Range("E32").Value = '  +3.20%  '
$ws.Range("D33").Value = '3.64'
$ws.Range("E33").Value = '  +4.65%  '
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("D35").Value = '1.447.18'
$ws.Range("D36").Value = '0.655'
$ws.Range("E36").Value = '  +1.88%  '
$ws.Range("E37").Value = '  +9.99%  '
$ws.Range("E38").Value = '  +4.33%  '
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("D40").Value = '80.21'
$ws.Range("E40").Value = '  +3.36%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").Value = '0.924'
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").Value = '13.51'
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = '6.07'
$ws.Range("E45").Value = '  +4.69%  '
$ws.Range("E46").Value = '  +1.88%  '
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("D49").Value = '1.944.75'
$ws.Range("E49").Value = '  +1.46%  '
$ws.Range("D50").Value = '106.22'
$ws.Range("E51").Value = '  +0.13%  '
